# Populate the values first.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the "header cell" look (bold font, centered/top-aligned, thin box
# border) on B1 first - this yields exactly one new font/border/cellXf.
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1").VerticalAlignment = -4160     # xlTop
$ws.Range("B1").Borders.LineStyle = 1         # xlContinuous
$ws.Range("B1").Borders.Weight = 2            # xlThin

# Re-use that exact style for A2 via copy/paste-format instead of
# re-applying each property (which would otherwise synthesize a second,
# separate cellXf for the same visual style).
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)           # xlPasteFormats
$excel.CutCopyMode = $false
